$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 5.624
$ws.Range("D2").Value = -7.56
$ws.Range("A3").Value = -21.535
$ws.Range("C3").Value = -12.505
$ws.Range("D6").Value = -7.855
$ws.Range("E8").Value = 16.636
$ws.Range("E9").Value = 16.461
$ws.Range("C12").Value = -11.536
$ws.Range("A14").Value = -21.624
$ws.Range("A16").Value = -21.222
$ws.Range("B18").Value = 5.414
$ws.Range("D19").Value = -8.078999999999999
$ws.Range("A21").Value = -20.662
$ws.Range("A23").Value = -20.254
$ws.Range("E23").Value = 16.586
$ws.Range("B24").Value = 7.05
$ws.Range("C24").Value = -13.093
$ws.Range("D24").Value = -7.220000000000001
$ws.Range("A25").Value = -20.55
$ws.Range("B25").Value = 6.325
$ws.Range("C25").Value = -12.523
$ws.Range("A26").Value = -21.349
$ws.Range("E26").Value = 16.625
$ws.Range("B27").Value = 5.564
$ws.Range("D27").Value = -8.062000000000001
$ws.Range("A29").Value = -21.162
$ws.Range("B30").Value = 5.712000000000001
$ws.Range("D30").Value = -7.377
$ws.Range("B31").Value = 4.976000000000001
$ws.Range("D31").Value = -7.671000000000001
$ws.Range("D33").Value = -7.666999999999999
$ws.Range("E37").Value = 16.485
$ws.Range("B39").Value = 7.798
$ws.Range("A40").Value = -20.27
$ws.Range("C41").Value = -12.61
$ws.Range("B42").Value = 8.395
$ws.Range("D42").Value = -8.472
$ws.Range("B48").Value = 5.176
$ws.Range("E48").Value = 17.303
$ws.Range("C50").Value = -13.087
$ws.Range("B51").Value = 5.216
$ws.Range("B52").Value = 4.933
$ws.Range("A53").Value = -21.814
$ws.Range("C53").Value = -11.891
$ws.Range("E54").Value = 16.766
$ws.Range("B55").Value = 4.581999999999999
$ws.Range("D55").Value = -8.131000000000002
$ws.Range("B56").Value = 5.546000000000001
$ws.Range("C56").Value = -12.676
$ws.Range("A57").Value = -21.338
$ws.Range("B57").Value = 5.917
$ws.Range("C57").Value = -12.531
$ws.Range("C58").Value = -12.826
$ws.Range("D58").Value = -8.094999999999999
$ws.Range("A59").Value = -21.852
$ws.Range("B60").Value = 5.005000000000001
$ws.Range("C61").Value = -13.181
$ws.Range("E62").Value = 16.867
$ws.Range("C63").Value = -11.757
$ws.Range("C64").Value = -11.805
$ws.Range("A65").Value = -21.564
$ws.Range("D65").Value = -7.639
$ws.Range("E65").Value = 16.752
$ws.Range("E66").Value = 17.087
$ws.Range("A69").Value = -21.56
$ws.Range("C70").Value = -11.938
$ws.Range("D70").Value = -7.569999999999999
$ws.Range("C72").Value = -11.799
$ws.Range("B73").Value = 6.729000000000001
$ws.Range("B74").Value = 8.492000000000001
$ws.Range("D74").Value = -8.044999999999998
$ws.Range("D75").Value = -7.557999999999998
$ws.Range("E75").Value = 16.701
$ws.Range("A79").Value = -21.237
$ws.Range("E81").Value = 16.604
$ws.Range("A83").Value = -21.372
$ws.Range("D83").Value = -8.526
$ws.Range("D84").Value = -7.934
$ws.Range("C86").Value = -12.574
$ws.Range("D86").Value = -7.965000000000001
$ws.Range("B89").Value = 5.214
$ws.Range("C89").Value = -11.873
$ws.Range("E89").Value = 17.21
$ws.Range("B90").Value = 5.338000000000001
$ws.Range("A91").Value = -21.527
$ws.Range("B92").Value = 5.962999999999999
$ws.Range("A93").Value = -21.324
$ws.Range("E94").Value = 17.408
$ws.Range("D96").Value = -7.533999999999999
$ws.Range("D97").Value = -8.020999999999999
$ws.Range("C98").Value = -12.712
$ws.Range("A100").Value = -21.473
$ws.Range("C100").Value = -12.279
$ws.Range("C102").Value = -12.884
